$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new (empty) column before column N (14) - "Late" and the
# following columns shift one position to the right.
$ws.Columns.Item(14).Insert()
$ws.Columns.Item(14).ColumnWidth = 10.14

# Make "Repayment schedule" the active sheet/tab and select cell R7,
# mirroring the saved view state in the workbook.
$ws.Activate()
$ws.Range("R7").Select()
